# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on the zh-cn and
# de-de report sheets to reflect the freshly-generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-18 10:32:08"
$wsZhCn.Range("D3").Value = "2016-02-18 10:32:08"
$wsZhCn.Range("G2").Value = "2016-02-18 10:33:02"
$wsZhCn.Range("G3").Value = "2016-02-18 10:33:02"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-18 10:32:19"
$wsDeDe.Range("D3").Value = "2016-02-18 10:32:19"
$wsDeDe.Range("G2").Value = "2016-02-18 10:33:25"
$wsDeDe.Range("G3").Value = "2016-02-18 10:33:25"
